$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") for data rows 2 through 494 all change
# from serial date 45172 (2023-09-03) to 45175 (2023-09-06).
$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(45175)

for ($row = 2; $row -le 494; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
